$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: Registro de Usuario -> Cerrar Sesion / Logout ---
$ws.Range("C8").Value = "Cerrar Sesión / Logout"
$ws.Range("D8").Value = "Se puede cerrar sesión y se redirige a la pagina principal."
$ws.Range("E8").Value = "Se puede cerrar sesión y se redirige a la pagina principal."
$ws.Range("E8").WrapText = $true

# --- Row 9: Uso del Navbar -> Avatar ---
$ws.Range("C9").Value = "Avatar"
$ws.Range("D9").Value = "Se puede Editar y Actualizar el Avatar."
$ws.Range("E9").Value = "En la seccion de Editar, una vez logueado se puede agregar un Avatar o modificar uno existente"

# --- Row 10: Carga de Eventos -> Registro de Usuario ---
$ws.Range("C10").Value = "Registro de Usuario"
$ws.Range("D10").Value = "Abre el formulario de Registro. Se puede completar el formulario y redirige a la pagina principal. En caso contrario, se presiona el boton ""Cancelar"" y se cierra el formulario."
$ws.Range("E10").Value = "Si el Registro el satisfactorio se puede iniciar sesión"
$ws.Range("E10").WrapText = $false

# --- Row 11: Eliminar Eventos -> Uso del Navbar ---
$ws.Range("C11").Value = "Uso del Navbar"
$ws.Range("D11").Value = "Los link redirigen a las secciones de la misma pagina."
$ws.Range("E11").Value = "Todos van a una página o sección distinta"

# --- Row 12: Edicion de Instrumento -> CRUD de Comentarios ---
$ws.Range("C12").Value = "CRUD de Comentarios"
$ws.Range("D12").Value = "Se pueden modificar y/o actualizar los Comentarios"
$ws.Range("E12").Value = "Si el usuario esta logueado ve los botones para interactuar, sino directamente lee los comentarios"

# --- Row 13: Eliminacion de Instrumento Musical -> CRUD de Eventos ---
$ws.Range("C13").Value = "CRUD de Eventos"
$ws.Range("C13").WrapText = $false
$ws.Range("D13").Value = "Se puede crear nuevos eventos con titulo y una imagen. Tambien se pueden editar eventos existentes o directamente eliminarlos"
$ws.Range("E13").Value = "Desde la vista del editor, hay un boton en cada evento que permite eliminarlo, editarlo. La visualizacion es en formato de Galería"

# --- Row 14: Edicion de Perfil de Usuario -> CRUD de Indumentaria ---
$ws.Range("C14").Value = "CRUD de Indumentaria"
$ws.Range("D14").Value = "Se pueden crear nuevas cards de indumentaria con titulo, descripción y una imagen. Tambien se pueden editar cards de indumentaria existentes o directamente eliminarlas."
$ws.Range("E14").Value = "El usuario logueado ve en cada Cards los botones para interactuar y el visitante solo ve una Galeria de cards"

# --- Rows 15 & 16: clear contents (keep as blank formatted rows, no explicit height) ---
$ws.Range("A15:F16").Delete(-4162)
$ws.Range("A15:F16").Insert(-4121)
$ws.Range("C15").WrapText = $true

# --- Rows 17 & 18: new blank rows, formatted like row 16 ---
$ws.Range("A16:F16").Copy()
$ws.Range("A17:F18").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row heights ---
$ws.Rows.Item(7).RowHeight = 138
$ws.Rows.Item(8).RowHeight = 38.25
$ws.Rows.Item(9).RowHeight = 38.25
$ws.Rows.Item(10).RowHeight = 75
$ws.Rows.Item(12).RowHeight = 30
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(14).RowHeight = 75

# --- Selection ---
$ws.Range("F9").Select() | Out-Null
